# Update the "handback-status" report timestamps to reflect a new
# report-generation run (commit message: "Generate Report for Handback").
#
# Sheets: "Overview", "zh-cn", "de-de"
#
#   Overview!G4  (Latest HO Xliff Generate Date)      2016-08-29 08:47:52 -> 2016-08-29 08:48:53
#   zh-cn!H4     (Correspond Handoff Datetime)        2016-08-29 08:47:47 -> 2016-08-29 08:48:49
#   zh-cn!K4     (Correspond Handback DateTime)       2016-08-29 08:48:20 -> 2016-08-29 08:49:17
#   de-de!H4     (Correspond Handoff Datetime)        2016-08-29 08:47:52 -> 2016-08-29 08:48:53
#   de-de!K4     (Correspond Handback DateTime)       2016-08-29 08:48:27 -> 2016-08-29 08:49:24

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$wsOverview.Range("G4").Value = "2016-08-29 08:48:53"

$wsZhCn.Range("H4").Value = "2016-08-29 08:48:49"
$wsZhCn.Range("K4").Value = "2016-08-29 08:49:17"

$wsDeDe.Range("H4").Value = "2016-08-29 08:48:53"
$wsDeDe.Range("K4").Value = "2016-08-29 08:49:24"
